# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the upstream data source.

$wb = $excel.ActiveWorkbook

# Map of row -> new F-column value for sheet "展览"
$zhanlanUpdates = @{
    2  = 15018
    3  = 19065
    5  = 138
    14 = 162
    15 = 219
    16 = 66
    17 = 1470
    22 = 7947
    29 = 6057
    30 = 116
    31 = 73
    34 = 284
    35 = 5434
    36 = 291
    37 = 12
    39 = 47
}

# Map of row -> new F-column value for sheet "全部类型"
$quanbuUpdates = @{
    2  = 15018
    3  = 19065
    5  = 138
    14 = 162
    15 = 219
    16 = 66
    17 = 1470
    23 = 7947
    32 = 6057
    33 = 116
    34 = 73
    37 = 284
    38 = 5434
    39 = 291
    40 = 12
    42 = 47
}

$wsZhanlan = $wb.Worksheets.Item("展览")
foreach ($row in $zhanlanUpdates.Keys) {
    $wsZhanlan.Range("F$row").Value = $zhanlanUpdates[$row]
}

$wsQuanbu = $wb.Worksheets.Item("全部类型")
foreach ($row in $quanbuUpdates.Keys) {
    $wsQuanbu.Range("F$row").Value = $quanbuUpdates[$row]
}
